# Update "horarios" (schedules) workbook with the latest scrape results.
# New scrape timestamp
$nuevaHora = "01:16:06"

$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 --------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $nuevaHora"

# Row 6 (first upcoming arrival) now corresponds to 14_ABASTO
$ws1.Cells.Item(6, 1).Value = $nuevaHora
$ws1.Cells.Item(6, 2).Value = "01:58"
$ws1.Cells.Item(6, 3).Value = "14_ABASTO"
$ws1.Cells.Item(6, 4).Value = 42

# Row 7 (second upcoming arrival) now corresponds to 215_ALUAR
$ws1.Cells.Item(7, 1).Value = $nuevaHora
$ws1.Cells.Item(7, 2).Value = "02:58"
$ws1.Cells.Item(7, 3).Value = "215_ALUAR"
$ws1.Cells.Item(7, 4).Value = 102

# --- Sheet 2: LP1912-215 ----------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $nuevaHora"

$ws2.Cells.Item(6, 1).Value = $nuevaHora
$ws2.Cells.Item(6, 2).Value = "02:58"
$ws2.Cells.Item(6, 4).Value = 102

# --- Sheet 3: 6203-6173 ------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $nuevaHora"
